$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "clave"

$ws.Range("A3").Value = "admin"
$ws.Range("B3").Value = "donpedro"

$ws.Range("A4").Value = "javier"
$ws.Range("B4").Value = 1234

$ws.Range("A5").Value = "tato"
$ws.Range("B5").Value = 1234

$ws.Range("A1:B1").ClearFormats()

$ws.Range("B5").Select()
